# Add a new "time_taken" column (F) with per-row timestamps to the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - same style as the other header cells (B1:E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Per-row timestamps (stored as text, matching the source data's inline strings)
$timestamps = @(
    "2021-10-05 13:39:35.221487",
    "2021-10-05 13:39:35.221498",
    "2021-10-05 13:39:35.221502",
    "2021-10-05 13:39:35.221505",
    "2021-10-05 13:39:35.221508",
    "2021-10-05 13:39:35.221511",
    "2021-10-05 13:39:35.221513",
    "2021-10-05 13:39:35.221515",
    "2021-10-05 13:39:35.221518",
    "2021-10-05 13:39:35.221521",
    "2021-10-05 13:39:35.221523",
    "2021-10-05 13:39:35.221526",
    "2021-10-05 13:39:35.221529",
    "2021-10-05 13:39:35.221531"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
